$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.270.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.560.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.68%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.16%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.476"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0611"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.243"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.80"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.777.23"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.566.70"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.99"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.506"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.280.65"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "59.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.31%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0711"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "186.63"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.85"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.85"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.66"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.91"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.39"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0465"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.081.22"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.32"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.495"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.766"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -9.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.795"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "92.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.691.60"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0504"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.404"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.06%  "
